# Generate Report for Handoff
#
# Refreshes the handoff timestamps for the source file
# c1b996f8-b300-4010-a124-dd7686a1d3e0.md (row 4 on every sheet):
#   - Overview!G4 "Latest HO Xliff Generate Date"
#   - zh-cn!H4    "Latest Handoff Datetime"
#   - de-de!H4    "Latest Handoff Datetime"

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$overview.Range("G4").Value = "2016-10-27 05:50:12"
$zhcn.Range("H4").Value = "2016-10-27 05:49:59"
$dede.Range("H4").Value = "2016-10-27 05:50:12"
